$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = 298
$ws.Range("D6").Value = 270
$ws.Range("E6").Value = 28
$ws.Range("F6").Value = 59.08096280087527
$ws.Range("G6").Value = 9.395973154362416
$ws.Range("H6").Value = 90.60402684563759
